$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Locate the anchor paragraph: the "Model 2" Heading1 paragraph in
#    the SAS "Model 2" code section (second "Model 2" occurrence in the
#    document - the first is "Model 2." used earlier as a report
#    heading, this one introduces the SAS/logistic-regression code
#    block that we are editing).
# ------------------------------------------------------------------
$headingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text.TrimEnd([char]13) -eq "Model 2" -and $pp.Style.NameLocal -eq "Heading 1") {
        $headingPara = $pp
    }
}
if ($headingPara -eq $null) {
    Write-Output "ERROR: could not find Model 2 heading"
}

$headingIndex = $headingPara.Range.Information(3)  # wdActiveEndAdjustedPageNumber not useful; use index search below instead

# Find the paragraph index numerically (Paragraphs collection has no direct IndexOf,
# so re-derive it by scanning again and comparing Range.Start).
$targetStart = $headingPara.Range.Start
$idx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Start -eq $targetStart) {
        $idx = $i
        break
    }
}
Write-Output ("Model 2 (code) heading paragraph index: " + $idx)

# The paragraph immediately before the heading is the empty paragraph
# that currently carries the _GoBack bookmark.
$prevPara = $d.Paragraphs.Item($idx - 1)
Write-Output ("Prev para text: >" + $prevPara.Range.Text + "<")

# The paragraph right after the heading is an empty <w:p/> that must be
# removed (the FILENAME line should follow the heading directly).
$afterHeading = $d.Paragraphs.Item($idx + 1)
Write-Output ("After-heading para text: >" + $afterHeading.Range.Text + "<")

# ------------------------------------------------------------------
# 2. Remove the _GoBack bookmark from its current (empty) location.
#    It will be re-inserted later, right after the updated model
#    statement text.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
    Write-Output "Removed old _GoBack bookmark"
}

# ------------------------------------------------------------------
# 3. Delete the empty paragraph right after the "Model 2" heading
#    (between the heading and the FILENAME line).
# ------------------------------------------------------------------
$afterHeading = $d.Paragraphs.Item($idx + 1)
$afterHeading.Range.Delete()
Write-Output "Deleted empty paragraph after heading"

# ------------------------------------------------------------------
# 4. Delete the empty paragraph between the FILENAME line and the
#    PROC IMPORT line.
# ------------------------------------------------------------------
$filenamePara = $d.Paragraphs.Item($idx + 1)
Write-Output ("FILENAME para: >" + $filenamePara.Range.Text + "<")
$afterFilename = $d.Paragraphs.Item($idx + 2)
Write-Output ("After FILENAME para: >" + $afterFilename.Range.Text + "<")
$afterFilename.Range.Delete()
Write-Output "Deleted empty paragraph after FILENAME line"

# ------------------------------------------------------------------
# 5. Fix the "drop" -> "DELETE" typo in the IF statement.
# ------------------------------------------------------------------
$d.Content.Find.Execute('IF shot_made_flag = "." THEN drop;', $true, $false, $false, $false, $false, $true, 1, $false, 'IF shot_made_flag = "." THEN DELETE;', 2) | Out-Null
Write-Output "Replaced drop -> DELETE"

# ------------------------------------------------------------------
# 6. Insert a new empty paragraph between "proc logistic data = Kobe;"
#    and the "model ..." statement paragraph.
# ------------------------------------------------------------------
$found = $d.Content.Find.Execute("proc logistic data = Kobe;")
if ($found) {
    $procPara = $d.Range($d.Content.Find.Parent.Start, $d.Content.Find.Parent.Start).Paragraphs.Item(1)
}
# Locate it robustly by scanning paragraphs text.
$procIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text.TrimEnd([char]13) -eq "proc logistic data = Kobe;") {
        $procIdx = $i
        break
    }
}
Write-Output ("proc logistic paragraph index: " + $procIdx)
$procPara = $d.Paragraphs.Item($procIdx)
$procPara.Range.InsertParagraphAfter()
Write-Output "Inserted empty paragraph after proc logistic line"

# ------------------------------------------------------------------
# 7. Update the "model" statement paragraph text and re-insert the
#    _GoBack bookmark right after the updated prefix.
# ------------------------------------------------------------------
$modelIdx = $procIdx + 2
$modelPara = $d.Paragraphs.Item($modelIdx)
Write-Output ("Model paragraph text before: >" + $modelPara.Range.Text + "<")

# a) drop the stray "= " right after "model"
$d.Content.Find.Execute("model = shot_made_flag (event = `'1`') =  shot_distance;", $true, $false, $false, $false, $false, $true, 1, $false, "`"Kobe Shots`": model shot_made_flag (event = `'1`') = shot_distance;", 2) | Out-Null

Write-Output ("Model paragraph text after: >" + $modelPara.Range.Text + "<")

# b) insert the _GoBack bookmark right after `(event = '1') ` (i.e. right
#    before "= shot_distance;")
$markerText = "`"Kobe Shots`": model shot_made_flag (event = `'1`') "
$found2 = $d.Content.Find.Execute($markerText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $bmRange = $d.Range($d.Content.Find.Parent.End, $d.Content.Find.Parent.End)
}
# Robust approach: locate via paragraph range + text offset.
$modelPara = $d.Paragraphs.Item($modelIdx)
$pStart = $modelPara.Range.Start
$insertPos = $pStart + $markerText.Length
$bmRange = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
Write-Output "Re-added _GoBack bookmark"

Write-Output ("Final model paragraph text: >" + $d.Paragraphs.Item($modelIdx).Range.Text + "<")
Write-Output ("Total paragraphs: " + $d.Paragraphs.Count)
